$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Footer containing the Pearson Edexcel logo with docPr id="3" (physical footer1.xml)
# -> rename the inline picture from image1.png to image2.png
$ftrA = $sec.Footers.Item(2)
$picA = $ftrA.Range.InlineShapes.Item(1)
$picA.Name = "image2.png"

# Footer containing the Pearson Edexcel logo with docPr id="2" (physical footer2.xml)
# -> rename the inline picture from image1.png to image2.png
$ftrB = $sec.Footers.Item(1)
$picB = $ftrB.Range.InlineShapes.Item(1)
$picB.Name = "image2.png"

# Header containing the BTec logo (docPr id="1")
# -> rename the inline picture from image2.jpg to image1.jpg
$hdr = $sec.Headers.Item(2)
$picC = $hdr.Range.InlineShapes.Item(1)
$picC.Name = "image1.jpg"
